$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new daily-log entry as row 17 (A:E), mirroring the existing rows
$ws.Range("A17").Value = 6
$ws.Range("B17").Value = "2：04-6;07"
$ws.Range("C17").Value = "函数 传值函数 传址函数 函数的嵌套调用和链式访问"
$ws.Range("E17").Value = "（这两天好像更多是概念和复习，新代码敲得少，明天应该好点，函数的递归了）"

# Match the author's final selection state on the new row
$ws.Range("E17").Select()
